$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data for FiftyForward Madison Station Senior Center
$ws.Range("A19").Value = "FiftyForward Madison Station Senior Center"
$ws.Range("B19").Value = "Davidson"
$ws.Range("C19").Value = "https://reports.mysidewalk.com/b797651eb3"

# Add hyperlink for the new sharing link, then restore the usual
# "Hyperlink" cell style used by the rest of column C
$ws.Hyperlinks.Add($ws.Range("C19"), "https://reports.mysidewalk.com/b797651eb3")
$ws.Range("C19").Style = $ws.Range("C18").Style

# Match the new selection left behind in the saved workbook
$ws.Range("A13").Select()
